# Swap the visit data (columns G through CG) between row 7 and row 8.
# Columns A-F (sampleid, program, location, county, state, full_state) and
# CH-CP (program_county/lat/lon/etc.) are per-sample identity/location fields
# and are not swapped - only the visit-specific data columns G:CG are.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Visit-data columns run from G (7) through CG (85); identity/location
# columns A:F and CH:CP stay put for both rows.
$firstCol = 7
$lastCol = 85
$row1 = 7
$row2 = 8

for ($col = $firstCol; $col -le $lastCol; $col++) {
    $cell1 = $ws.Cells.Item($row1, $col)
    $cell2 = $ws.Cells.Item($row2, $col)

    $val1 = $cell1.Value2
    $val2 = $cell2.Value2

    # Only touch cells whose content actually differs between the two rows;
    # this avoids clobbering identical blank/empty-string cells (assigning an
    # empty string would clear the cell instead of leaving it untouched).
    if ($val1 -ne $val2) {
        $cell1.Value2 = $val2
        $cell2.Value2 = $val1
    }
}
